$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ICtPSFfL")

# --- New rows 9-11: heavy/residual fuel oil, LPG propane/butane, hydrogen ---
# Each gets a label in column A and "=0" placeholder formulas across B:AK (years 2015-2050),
# matching the pattern already used by existing rows (e.g. row 8 "jet fuel").
$newRows = @(
    @{ Row = 9;  Label = "heavy or residual fuel oil" },
    @{ Row = 10; Label = "LPG propane or butane" },
    @{ Row = 11; Label = "hydrogen" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Label
    for ($c = 2; $c -le 37; $c++) {
        $ws.Cells.Item($r, $c).Formula = "=0"
    }
}

# --- A1: new header label "Cost ($/BTU)" in bold (matches style used elsewhere, e.g. Calcs!A1) ---
$ws.Range("A1").Value = "Cost (`$/BTU)"
$ws.Range("A1").Font.Bold = $true

# --- Column A width grows to fit the new label ---
$ws.Columns.Item(1).ColumnWidth = 23
